$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Sheet view: scroll/selection moved while entering the remaining 2016 (run14) embryo rows ---
$win = $excel.ActiveWindow
$win.ScrollRow = 154
$win.ScrollColumn = 1
$ws.Range("L233").Select()

# --- Update Side1 (E) values where they changed, and fill in Yolk/Body Area/Count (F-I) for 2016 run14 embryos ---

# Row 203
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 353562
$arr[0,1] = 48
$arr[0,2] = 523597
$arr[0,3] = 199
$ws.Range("F203:I203").Value = $arr

# Row 204
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 313698
$arr[0,1] = 150
$arr[0,2] = 629227
$arr[0,3] = 287
$ws.Range("F204:I204").Value = $arr

# Row 205
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 120790
$arr[0,1] = 55
$arr[0,2] = 656212
$arr[0,3] = 219
$ws.Range("F205:I205").Value = $arr

# Row 206, Side1 corrected to L
$ws.Range("E206").Value = "L"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 148627
$arr[0,1] = 50
$arr[0,2] = 755045
$arr[0,3] = 305
$ws.Range("F206:I206").Value = $arr

# Row 207
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 306212
$arr[0,1] = 83
$arr[0,2] = 451840
$arr[0,3] = 184
$ws.Range("F207:I207").Value = $arr

# Row 208
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 394277
$arr[0,1] = 90
$arr[0,2] = 652741
$arr[0,3] = 301
$ws.Range("F208:I208").Value = $arr

# Row 209
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 214232
$arr[0,1] = 74
$arr[0,2] = 10224290
$arr[0,3] = 335
$ws.Range("F209:I209").Value = $arr

# Row 210
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 259206
$arr[0,1] = 44
$arr[0,2] = 534428
$arr[0,3] = 119
$ws.Range("F210:I210").Value = $arr

# Row 211
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 176261
$arr[0,1] = 137
$arr[0,2] = 756448
$arr[0,3] = 325
$ws.Range("F211:I211").Value = $arr

# Row 212
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 184731
$arr[0,1] = 150
$arr[0,2] = 495103
$arr[0,3] = 309
$ws.Range("F212:I212").Value = $arr

# Row 213
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 183509
$arr[0,1] = 35
$arr[0,2] = 331943
$arr[0,3] = 14
$ws.Range("F213:I213").Value = $arr

# Row 214
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 139952
$arr[0,1] = 49
$arr[0,2] = 593964
$arr[0,3] = 311
$ws.Range("F214:I214").Value = $arr

# Row 215, Side1 corrected to L
$ws.Range("E215").Value = "L"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 189762
$arr[0,1] = 64
$arr[0,2] = 606373
$arr[0,3] = 197
$ws.Range("F215:I215").Value = $arr

# Row 216
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 273636
$arr[0,1] = 64
$arr[0,2] = 460342
$arr[0,3] = 151
$ws.Range("F216:I216").Value = $arr

# Row 217
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 144104
$arr[0,1] = 110
$arr[0,2] = 616080
$arr[0,3] = 250
$ws.Range("F217:I217").Value = $arr

# Row 218
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 139398
$arr[0,1] = 29
$arr[0,2] = 549610
$arr[0,3] = 176
$ws.Range("F218:I218").Value = $arr

# Row 219, Side1 corrected to R
$ws.Range("E219").Value = "R"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 264966
$arr[0,1] = 56
$arr[0,2] = 618484
$arr[0,3] = 173
$ws.Range("F219:I219").Value = $arr

# Row 220, Side1 corrected to R
$ws.Range("E220").Value = "R"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 167438
$arr[0,1] = 128
$arr[0,2] = 611786
$arr[0,3] = 263
$ws.Range("F220:I220").Value = $arr

# Row 221, Side1 corrected to L
$ws.Range("E221").Value = "L"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 166458
$arr[0,1] = 108
$arr[0,2] = 788662
$arr[0,3] = 229
$ws.Range("F221:I221").Value = $arr

# Row 222
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 286996
$arr[0,1] = 84
$arr[0,2] = 496598
$arr[0,3] = 116
$ws.Range("F222:I222").Value = $arr

# Row 223, Side1 corrected to R
$ws.Range("E223").Value = "R"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 152683
$arr[0,1] = 88
$arr[0,2] = 619601
$arr[0,3] = 240
$ws.Range("F223:I223").Value = $arr

# Row 224
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 226233
$arr[0,1] = 104
$arr[0,2] = 603649
$arr[0,3] = 304
$ws.Range("F224:I224").Value = $arr

# Row 225
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 141389
$arr[0,1] = 42
$arr[0,2] = 705803
$arr[0,3] = 247
$ws.Range("F225:I225").Value = $arr

# Row 226: NA embryo - mirror Side2 NA columns (F:I) and add NA placeholders (P:S)
$ws.Range("F226:I226").Value = "NA"
$ws.Range("P226:S226").Value = "NA"

# Row 227
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 267069
$arr[0,1] = 128
$arr[0,2] = 505954
$arr[0,3] = 229
$ws.Range("F227:I227").Value = $arr

# Row 228, Side1 corrected to R
$ws.Range("E228").Value = "R"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 223385
$arr[0,1] = 63
$arr[0,2] = 819133
$arr[0,3] = 161
$ws.Range("F228:I228").Value = $arr

# Row 229
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 101321
$arr[0,1] = 61
$arr[0,2] = 595781
$arr[0,3] = 341
$ws.Range("F229:I229").Value = $arr

# Row 230
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 187463
$arr[0,1] = 91
$arr[0,2] = 576577
$arr[0,3] = 144
$ws.Range("F230:I230").Value = $arr

# Row 231
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 226767
$arr[0,1] = 63
$arr[0,2] = 556146
$arr[0,3] = 219
$ws.Range("F231:I231").Value = $arr

# Row 232: NA embryo - mirror Side2 NA columns (F:I) and add NA placeholders (P:S)
$ws.Range("F232:I232").Value = "NA"
$ws.Range("P232:S232").Value = "NA"

# Row 233, Side1 corrected to R
$ws.Range("E233").Value = "R"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 233934
$arr[0,1] = 33
$arr[0,2] = 439524
$arr[0,3] = 133
$ws.Range("F233:I233").Value = $arr

# Row 234
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 170029
$arr[0,1] = 81
$arr[0,2] = 500161
$arr[0,3] = 289
$ws.Range("F234:I234").Value = $arr

# Row 235
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 167129
$arr[0,1] = 80
$arr[0,2] = 543545
$arr[0,3] = 305
$ws.Range("F235:I235").Value = $arr

# Row 236, Side1 corrected to R
$ws.Range("E236").Value = "R"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 132168
$arr[0,1] = 110
$arr[0,2] = 614676
$arr[0,3] = 359
$ws.Range("F236:I236").Value = $arr

# Row 237, Side1 corrected to L
$ws.Range("E237").Value = "L"
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 162133
$arr[0,1] = 100
$arr[0,2] = 739948
$arr[0,3] = 392
$ws.Range("F237:I237").Value = $arr

# Row 238
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 306476
$arr[0,1] = 143
$arr[0,2] = 441301
$arr[0,3] = 228
$ws.Range("F238:I238").Value = $arr

# Row 239
$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 322644
$arr[0,1] = 154
$arr[0,2] = 523646
$arr[0,3] = 258
$ws.Range("F239:I239").Value = $arr
